$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2559.1155
$ws.Range("J17").Value = 2409.48
$ws.Range("L17").Value = 7228.440000000001
$ws.Range("N17").Value = -7564.440000000001
$ws.Range("H28").Value = 692.2353000000001
$ws.Range("J28").Value = 3001.6667
$ws.Range("L28").Value = 3001.6667
$ws.Range("N28").Value = -3971.6667
$ws.Range("H62").Value = 7208.8
$ws.Range("I62").Value = 7208.8
$ws.Range("K62").Value = 7208.8
$ws.Range("M62").Value = -6584.8
$ws.Range("H65").Value = 7208.8
$ws.Range("I65").Value = 7208.8
$ws.Range("K65").Value = 36044
$ws.Range("M65").Value = -32924
$ws.Range("H98").Value = 1598.7916
$ws.Range("I98").Value = 1289.6364
$ws.Range("K98").Value = 1289.6364
$ws.Range("M98").Value = 208.3635999999999
$ws.Range("H112").Value = 2403.7646
$ws.Range("J112").Value = 2492.9375
$ws.Range("L112").Value = 7478.8125
$ws.Range("N112").Value = -9694.8125
$ws.Range("H122").Value = 1598.7916
$ws.Range("I122").Value = 1289.6364
$ws.Range("K122").Value = 3868.9092
$ws.Range("M122").Value = -1418.9092
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 1700.7407
$ws.Range("I132").Value = 1744.6086
$ws.Range("K132").Value = 5233.825800000001
$ws.Range("M132").Value = -2703.825800000001
$ws.Range("H137").Value = 9603.223
$ws.Range("J137").Value = 13551.632
$ws.Range("L137").Value = 40654.896
$ws.Range("N137").Value = -45754.896
$ws.Range("H138").Value = 2702.3333
$ws.Range("I138").Value = 2244.3572
$ws.Range("K138").Value = 6733.071599999999
$ws.Range("M138").Value = -1593.071599999999
$ws.Range("H141").Value = 3857.6
$ws.Range("I141").Value = 5649.5
$ws.Range("J141").Value = 2663
$ws.Range("K141").Value = 16948.5
$ws.Range("L141").Value = 7989
$ws.Range("M141").Value = -11768.5
$ws.Range("N141").Value = -18349

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -84
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 30000
$ws.Range("I6").Value = 30000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -29827
$ws.Range("N6").ClearContents()
$ws.Range("H61").Value = 6657.425
$ws.Range("I61").Value = 4808.909
$ws.Range("K61").Value = 4808.909
$ws.Range("M61").Value = -4596.909
$ws.Range("H76").Value = 36794
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 36794
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 36794
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -37470
$ws.Range("H79").Value = 36794
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 36794
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 36794
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -39134
$ws.Range("H97").Value = 1130.8462
$ws.Range("I97").Value = 1329.6842
$ws.Range("K97").Value = 1329.6842
$ws.Range("M97").Value = -833.6841999999999
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H136").Value = 6657.425
$ws.Range("I136").Value = 4808.909
$ws.Range("K136").Value = 14426.727
$ws.Range("M136").Value = -11876.727

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10838.219
$ws.Range("I20").Value = 13262
$ws.Range("J20").Value = 2181.8572
$ws.Range("K20").Value = 13262
$ws.Range("L20").Value = 2181.8572
$ws.Range("M20").Value = -13015
$ws.Range("N20").Value = -2675.8572
$ws.Range("H22").Value = 274.5
$ws.Range("I22").Value = 324
$ws.Range("J22").Value = 225
$ws.Range("K22").Value = 324
$ws.Range("L22").Value = 225
$ws.Range("M22").Value = -151
$ws.Range("N22").Value = -571

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 443
$ws.Range("J22").Value = 1150
$ws.Range("L22").Value = 1150
$ws.Range("N22").Value = -1850
$ws.Range("H31").Value = 3882.5144
$ws.Range("I31").Value = 3963.6365
$ws.Range("K31").Value = 3963.6365
$ws.Range("M31").Value = -3668.6365
$ws.Range("H34").Value = 3882.5144
$ws.Range("I34").Value = 3963.6365
$ws.Range("K34").Value = 3963.6365
$ws.Range("M34").Value = -3761.6365
$ws.Range("H106").Value = 49310.8
$ws.Range("J106").Value = 49310.8
$ws.Range("L106").Value = 49310.8
$ws.Range("N106").Value = -51834.8
$ws.Range("H141").Value = 466998.8
$ws.Range("J141").Value = 552500
$ws.Range("L141").Value = 552500
$ws.Range("N141").Value = -562860

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 46749132
$ws.Range("I4").Value = 48400884
$ws.Range("J4").Value = 500100
$ws.Range("K4").Value = 145202652
$ws.Range("L4").Value = 1500300
$ws.Range("M4").Value = -145202540
$ws.Range("N4").Value = -1500524
$ws.Range("H7").Value = 79628.57000000001
$ws.Range("I7").Value = 86.666664
$ws.Range("J7").Value = 139285
$ws.Range("K7").Value = 259.999992
$ws.Range("L7").Value = 417855
$ws.Range("M7").Value = -147.999992
$ws.Range("N7").Value = -418079
$ws.Range("H99").Value = 10631.25
$ws.Range("I99").Value = 5025
$ws.Range("K99").Value = 15075
$ws.Range("M99").Value = -12829
$ws.Range("H131").Value = 2606.8408
$ws.Range("I131").Value = 1177
$ws.Range("J131").Value = 2790.1538
$ws.Range("K131").Value = 3531
$ws.Range("L131").Value = 8370.4614
$ws.Range("M131").Value = 1509
$ws.Range("N131").Value = -18450.4614
$ws.Range("H138").Value = 2585
$ws.Range("I138").Value = 1377.5
$ws.Range("K138").Value = 4132.5
$ws.Range("M138").Value = 1007.5
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7423
$ws.Range("I70").Value = 8707.5
$ws.Range("J70").Value = 7056
$ws.Range("K70").Value = 8707.5
$ws.Range("L70").Value = 7056
$ws.Range("M70").Value = -8437.5
$ws.Range("N70").Value = -7596
$ws.Range("H73").Value = 7423
$ws.Range("I73").Value = 8707.5
$ws.Range("J73").Value = 7056
$ws.Range("K73").Value = 8707.5
$ws.Range("L73").Value = 7056
$ws.Range("M73").Value = -7771.5
$ws.Range("N73").Value = -8928
$ws.Range("H132").Value = 18555.428
$ws.Range("I132").Value = 15928.818
$ws.Range("J132").Value = 23000.46
$ws.Range("K132").Value = 47786.454
$ws.Range("L132").Value = 69001.38
$ws.Range("M132").Value = -45256.454
$ws.Range("N132").Value = -74061.38
$ws.Range("H141").Value = 63465.8
$ws.Range("J141").Value = 63465.8
$ws.Range("L141").Value = 63465.8
$ws.Range("N141").Value = -73825.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8499.277
$ws.Range("J7").Value = 9636.272000000001
$ws.Range("L7").Value = 9636.272000000001
$ws.Range("N7").Value = -9860.272000000001
$ws.Range("H46").Value = 1586.1613
$ws.Range("I46").Value = 980
$ws.Range("K46").Value = 980
$ws.Range("M46").Value = -792
$ws.Range("H93").Value = 3517.9546
$ws.Range("I93").Value = 4243.8125
$ws.Range("K93").Value = 4243.8125
$ws.Range("M93").Value = -2995.8125
$ws.Range("H126").Value = 8499.277
$ws.Range("J126").Value = 9636.272000000001
$ws.Range("L126").Value = 28908.816
$ws.Range("N126").Value = -33848.81600000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 6537836
$ws.Range("I107").Value = 1260.6364
$ws.Range("J107").Value = 18521558
$ws.Range("K107").Value = 3781.9092
$ws.Range("L107").Value = 55564674
$ws.Range("M107").Value = -1861.9092
$ws.Range("N107").Value = -55568514
$ws.Range("H122").Value = 5017.5557
$ws.Range("I122").Value = 3481.1304
$ws.Range("K122").Value = 10443.3912
$ws.Range("M122").Value = -7993.3912
$ws.Range("H124").Value = 44950
$ws.Range("J124").Value = 44950
$ws.Range("L124").Value = 44950
$ws.Range("N124").Value = -54770
$ws.Range("H126").Value = 7855.1387
$ws.Range("I126").Value = 5337.6895
$ws.Range("K126").Value = 16013.0685
$ws.Range("M126").Value = -13543.0685
$ws.Range("H132").Value = 177697.05
$ws.Range("I132").Value = 359331.06
$ws.Range("K132").Value = 1077993.18
$ws.Range("M132").Value = -1075463.18
$ws.Range("H136").Value = 3510790.2
$ws.Range("J136").Value = 2882.913
$ws.Range("L136").Value = 8648.739
$ws.Range("N136").Value = -13748.739
$ws.Range("H141").Value = 70000
$ws.Range("J141").Value = 70000
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360
